# Update the NATMI LR-pair sheet with the new TPM-based results.
# The underlying analysis now includes "ECs" as an additional sending/target
# cluster (in addition to the existing FAPs and MuSCs), so the 2x3 grid of
# Sending-cluster x Target-cluster rows grows into a full 3x3 grid (9 rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ligand / Receptor are constant across all rows.
$ligand   = "Slitrk1"
$receptor = "Ptprs"

# Each entry: SendingCluster, TargetCluster, followed by columns E..T
$rows = @(
    @("ECs",  "ECs",   2, 0.6666666666666666, 0.007204333333333333, 0.021613,
        0.01038998123717357, 0.01038998123717357, 3, 1,
        2.425633666666667, 7.276901000000001,
        0.0662600404061536, 0.06626004040615362,
        0.01747507347922222, 0.157275661313,
        0.0006884405765942984, 0.0006884405765942986),

    @("ECs",  "FAPs",  2, 0.6666666666666666, 0.007204333333333333, 0.021613,
        0.01038998123717357, 0.01038998123717357, 3, 1,
        15.503283, 46.509849,
        0.4234968256437875, 0.4234968256437876,
        0.111690818493, 1.005217366437,
        0.004400124072441518, 0.004400124072441519),

    @("ECs",  "MuSCs", 2, 0.6666666666666666, 0.007204333333333333, 0.021613,
        0.01038998123717357, 0.01038998123717357, 3, 1,
        18.67887366666667, 56.03662100000001,
        0.5102431339500588, 0.5102431339500588,
        0.1345688321858889, 1.211119489673,
        0.00530141658813775, 0.005301416588137751),

    @("FAPs", "ECs",   1, 0.3333333333333333, 0.008164333333333334, 0.024493,
        0.0117744788063708, 0.0117744788063708, 3, 1,
        2.425633666666667, 7.276901000000001,
        0.0662600404061536, 0.06626004040615362,
        0.01980368179922222, 0.178233136193,
        0.0007801774414715288, 0.000780177441471529),

    @("FAPs", "FAPs",  1, 0.3333333333333333, 0.008164333333333334, 0.024493,
        0.0117744788063708, 0.0117744788063708, 3, 1,
        15.503283, 46.509849,
        0.4234968256437875, 0.4234968256437876,
        0.126573970173, 1.139165731557,
        0.004986454398108089, 0.004986454398108089),

    @("FAPs", "MuSCs", 1, 0.3333333333333333, 0.008164333333333334, 0.024493,
        0.0117744788063708, 0.0117744788063708, 3, 1,
        18.67887366666667, 56.03662100000001,
        0.5102431339500588, 0.5102431339500588,
        0.1525005509058889, 1.372504958153,
        0.006007846966791188, 0.006007846966791188),

    @("MuSCs", "ECs",   3, 1, 0.6780236666666667, 2.034071,
        0.9778355399564556, 0.9778355399564557, 3, 1,
        2.425633666666667, 7.276901000000001,
        0.0662600404061536, 0.06626004040615362,
        1.644637032663445, 14.801733293971,
        0.06479142238808777, 0.06479142238808779),

    @("MuSCs", "FAPs",  3, 1, 0.6780236666666667, 2.034071,
        0.9778355399564556, 0.9778355399564557, 3, 1,
        15.503283, 46.509849,
        0.4234968256437875, 0.4234968256437876,
        10.511592785031, 94.604335065279,
        0.4141102471732379, 0.414110247173238),

    @("MuSCs", "MuSCs", 3, 1, 0.6780236666666667, 2.034071,
        0.9778355399564556, 0.9778355399564557, 3, 1,
        18.67887366666667, 56.03662100000001,
        0.5102431339500588, 0.5102431339500588,
        12.66471841267678, 113.982465714091,
        0.4989338703951299, 0.4989338703951299)
)

$rowIndex = 2
foreach ($r in $rows) {
    $ws.Cells.Item($rowIndex, 1).Value  = $r[0]        # A: Sending cluster
    $ws.Cells.Item($rowIndex, 2).Value  = $ligand       # B: Ligand symbol
    $ws.Cells.Item($rowIndex, 3).Value  = $receptor     # C: Receptor symbol
    $ws.Cells.Item($rowIndex, 4).Value  = $r[1]        # D: Target cluster

    for ($col = 5; $col -le 20; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $r[$col - 3]
    }

    $rowIndex++
}
